$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Collapse the "Honorific, First Name, Given Name, Surname, Lineage,
#    Other Name" run-sequence down to just "Lineage" (the styled runs for
#    Honorific / First Name / Given Name / Surname / Other Name - and the
#    ", " separator runs that went with them - are removed; the plain
#    "Lineage" run, the ", " run before it and the "." run after it stay).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Honorific, First Name, Given Name, Surname, ", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 2) | Out-Null

$d.Content.Find.Execute(
    ", Other Name", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the character styles that backed those runs - they are now
#    unused. Delete from the highest style index down to the lowest so
#    earlier deletions never invalidate an index this script still needs.
# ---------------------------------------------------------------------------

$unusedStyles = @(
    "Editor",
    "Author",
    "OtherName",
    "a",
    "Surname",
    "GivenName",
    "FirstName",
    "Honorific"
)

foreach ($styleName in $unusedStyles) {
    $d.Styles($styleName).Delete()
}
